$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("绩效表")

$ws.Range("K7").Value = 0.063
$ws.Range("K8").Value = 0.063
$ws.Range("K10").Value = 0.063
$ws.Range("K11").Value = 0.063
$ws.Range("K12").Value = 0.063
$ws.Range("K13").Value = 0.063
$ws.Range("K14").Value = 0.063
$ws.Range("K15").Value = 0.063
$ws.Range("K16").Value = 0.063
$ws.Range("K18").Value = 0.063

$ws.Range("J29").Value = 1.13
$ws.Range("H30").Value = 0.63
$ws.Range("I30").Value = "0.063+0.063+0.063+0.063+0.063+0.063+0.063+0.063+0.063+0.063=0.63"
